$d = $word.ActiveDocument

# 1. Fix the typo "pratices" -> "practices" in the workshop title.
$d.Content.Find.Execute("Best pratices in age", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Best practices in age", 2) | Out-Null

# 2. Reorder the author list so "Peter Comeau, " appears right after
#    "Daniel Ricard, " (i.e. before "Aaron Adamack, ") instead of after
#    "Jacob Burbank, ".
$peterRng = $d.Content
$peterRng.Find.Execute("Peter Comeau, ") | Out-Null
$peterRng.Delete()

$aaronRng = $d.Content
$aaronRng.Find.Execute("Aaron Adamack, ") | Out-Null
$aaronStart = $d.Range($aaronRng.Start, $aaronRng.Start)
$aaronStart.InsertBefore("Peter Comeau, ")
